$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook
$tmp = $wb.Worksheets.Add()

$tmp.Range("A1").Value = '!!!ObjTables schema=''SBtab'' objTablesVersion=''0.0.9'' date=''2020-04-26 21:09:37'''
$tmp.Range("A2").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Compartment'' name=''Compartment'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A3").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Compound'' name=''Compound'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A4").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Definition'' name=''Definition'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A5").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Enzyme'' name=''Enzyme'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A6").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A7").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Gene'' name=''Gene'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A8").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Layout'' name=''Layout'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A9").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Measurement'' name=''Measurement'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A10").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''PbConfig'' name=''PbConfig'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A11").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Position'' name=''Position'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A12").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Protein'' name=''Protein'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A13").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Quantity'' name=''Quantity'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A14").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A15").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A16").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Reaction'' name=''Reaction'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A17").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A18").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Regulator'' name=''Regulator'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A19").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Relation'' name=''Relation'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A20").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Relationship'' name=''Relationship'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A21").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A22").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A23").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A24").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A25").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A26").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A27").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-04-26 21:09:37'' objTablesVersion=''0.0.9'''
$tmp.Range("A28").Value = '!FromObject'
$tmp.Range("A29").Value = '!ToObject'

# --- Top-level ObjTables header (sheet1 !!Compartment, cell A1) ---
$tmp.Range("A1").Copy()
$wb.Worksheets.Item("!!Compartment").Range("A1").PasteSpecial(-4163)

# --- Per-table ObjTables header rows ---
$tmp.Range("A2").Copy()
$wb.Worksheets.Item("!!Compartment").Range("A2").PasteSpecial(-4163)
$tmp.Range("A3").Copy()
$wb.Worksheets.Item("!!Compound").Range("A1").PasteSpecial(-4163)
$tmp.Range("A4").Copy()
$wb.Worksheets.Item("!!Definition").Range("A1").PasteSpecial(-4163)
$tmp.Range("A5").Copy()
$wb.Worksheets.Item("!!Enzyme").Range("A1").PasteSpecial(-4163)
$tmp.Range("A6").Copy()
$wb.Worksheets.Item("!!FbcObjective").Range("A1").PasteSpecial(-4163)
$tmp.Range("A7").Copy()
$wb.Worksheets.Item("!!Gene").Range("A1").PasteSpecial(-4163)
$tmp.Range("A8").Copy()
$wb.Worksheets.Item("!!Layout").Range("A1").PasteSpecial(-4163)
$tmp.Range("A9").Copy()
$wb.Worksheets.Item("!!Measurement").Range("A1").PasteSpecial(-4163)
$tmp.Range("A10").Copy()
$wb.Worksheets.Item("!!PbConfig").Range("A1").PasteSpecial(-4163)
$tmp.Range("A11").Copy()
$wb.Worksheets.Item("!!Position").Range("A1").PasteSpecial(-4163)
$tmp.Range("A12").Copy()
$wb.Worksheets.Item("!!Protein").Range("A1").PasteSpecial(-4163)
$tmp.Range("A13").Copy()
$wb.Worksheets.Item("!!Quantity").Range("A1").PasteSpecial(-4163)
$tmp.Range("A14").Copy()
$wb.Worksheets.Item("!!QuantityInfo").Range("A1").PasteSpecial(-4163)
$tmp.Range("A15").Copy()
$wb.Worksheets.Item("!!QuantityMatrix").Range("A1").PasteSpecial(-4163)
$tmp.Range("A16").Copy()
$wb.Worksheets.Item("!!Reaction").Range("A1").PasteSpecial(-4163)
$tmp.Range("A17").Copy()
$wb.Worksheets.Item("!!ReactionStoichiometry").Range("A1").PasteSpecial(-4163)
$tmp.Range("A18").Copy()
$wb.Worksheets.Item("!!Regulator").Range("A1").PasteSpecial(-4163)
$tmp.Range("A19").Copy()
$wb.Worksheets.Item("!!Relation").Range("A1").PasteSpecial(-4163)
$tmp.Range("A20").Copy()
$wb.Worksheets.Item("!!Relationship").Range("A1").PasteSpecial(-4163)
$tmp.Range("A21").Copy()
$wb.Worksheets.Item("!!SparseMatrix").Range("A1").PasteSpecial(-4163)
$tmp.Range("A22").Copy()
$wb.Worksheets.Item("!!SparseMatrixColumn").Range("A1").PasteSpecial(-4163)
$tmp.Range("A23").Copy()
$wb.Worksheets.Item("!!SparseMatrixOrdered").Range("A1").PasteSpecial(-4163)
$tmp.Range("A24").Copy()
$wb.Worksheets.Item("!!SparseMatrixRow").Range("A1").PasteSpecial(-4163)
$tmp.Range("A25").Copy()
$wb.Worksheets.Item("!!StoichiometricMatrix").Range("A1").PasteSpecial(-4163)
$tmp.Range("A26").Copy()
$wb.Worksheets.Item("!!rxnconContingencyList").Range("A1").PasteSpecial(-4163)
$tmp.Range("A27").Copy()
$wb.Worksheets.Item("!!rxnconReactionList").Range("A1").PasteSpecial(-4163)

# --- !From -> !FromObject / !To -> !ToObject header cells ---
$tmp.Range("A28").Copy()
$wb.Worksheets.Item("!!Relation").Range("G2").PasteSpecial(-4163)
$tmp.Range("A28").Copy()
$wb.Worksheets.Item("!!Relationship").Range("B2").PasteSpecial(-4163)
$tmp.Range("A29").Copy()
$wb.Worksheets.Item("!!Relation").Range("H2").PasteSpecial(-4163)
$tmp.Range("A29").Copy()
$wb.Worksheets.Item("!!Relationship").Range("C2").PasteSpecial(-4163)

# --- cleanup staging sheet ---
[void]$tmp.Delete()

# --- Data validation title updates on !!Relation sheet ---
$wsRelation = $wb.Worksheets.Item("!!Relation")
$dvFrom = $wsRelation.Range("G2:G3").Validation
$dvFrom.ErrorTitle = "FromObject"
$dvFrom.InputTitle = "FromObject"
$dvTo = $wsRelation.Range("H2:H3").Validation
$dvTo.ErrorTitle = "ToObject"
$dvTo.InputTitle = "ToObject"

# --- Data validation title updates on !!Relationship sheet ---
$wsRelationship = $wb.Worksheets.Item("!!Relationship")
$dvFrom2 = $wsRelationship.Range("B2:B3").Validation
$dvFrom2.ErrorTitle = "FromObject"
$dvFrom2.InputTitle = "FromObject"
$dvTo2 = $wsRelationship.Range("C2:C3").Validation
$dvTo2.ErrorTitle = "ToObject"
$dvTo2.InputTitle = "ToObject"
